$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Measurements")

# Renumber the "Id" column (B) down by one: row 4 (was Id=1) becomes Id=0,
# row 5 (was Id=2) becomes Id=1, ... row 178 (was Id=175) becomes Id=174.
# All other cell content in these rows is left untouched.
for ($r = 4; $r -le 178; $r++) {
    $ws.Cells.Item($r, 2).Value = $r - 4
}

# Add the new trailing row (180th row overall / row 179) that now holds the
# Id=175 entry that got displaced by the renumbering above. Copy the
# formatting from the row directly above (row 178) so the new cells pick up
# the same styles (s="4" for the Id cell, s="28" for the M column), matching
# the pattern used throughout the rest of the table's tail rows.
$ws.Range("B178").Copy() | Out-Null
$ws.Range("B179").PasteSpecial(-4122) | Out-Null
$ws.Range("M178").Copy() | Out-Null
$ws.Range("M179").PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = 0

$ws.Cells.Item(179, 2).Value = 175

# Grow the "Table2" structured table by one row so it covers the new row.
$lo = $ws.ListObjects.Item(1)
$lo.Resize($ws.Range("B3:P179")) | Out-Null

# Restore the active selection (the commit moved it from R15 to C4).
$ws.Range("C4").Select() | Out-Null
